# Update to latest input spreadsheet
# Mirrors the upstream commit that refreshed the ITR Tool sample data:
#  - "ITR input data": NIPPON STEEL row (31) re-denominated from JPY to
#    USD (currency label + market cap/revenue/ev/evic/assets formulas),
#    a region fix on row 37, and newly supplied market_cap figures for
#    rows 45-47 (highlighted in red, matching the workbook's convention
#    for manually-entered data).
#  - View/selection bookkeeping on "ITR input data" and
#    "ITR target input data".
#  - "Portfolio" sheet's RANDBETWEEN-driven sample investment values
#    recalculate naturally (volatile formulas) once the workbook is
#    touched/recalculated.

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ITR input data")
$wsTarget = $wb.Worksheets.Item("ITR target input data")

# ---------------------------------------------------------------------
# ITR input data
# ---------------------------------------------------------------------
$wsInput.Activate()

# Row 31 (NIPPON STEEL CORP): currency changes from a stray "megaJPY"
# label to "USD", and the JPY-denominated figures are converted to USD
# (divide by the 107.92 JPY/USD rate and scale to actual units).
$wsInput.Range("H31").Value = "USD"
$wsInput.Range("J31").Formula = "=879400*1000000/107.92"
$wsInput.Range("K31").Formula = "=5921500*1000000/107.92"
$wsInput.Range("L31").Formula = "=M31+289549*1000000/107.92"
$wsInput.Range("M31").Formula = "=J31+2488741*1000000/107.92"
$wsInput.Range("N31").Formula = "=7444965*1000000/107.92"

# Row 37: region corrected to "Asia".
$wsInput.Range("E37").Value = "Asia"

# Rows 45-47: market_cap (column J) was missing; fill in with the newly
# sourced values, flagged in red like other manually-keyed-in entries.
$wsInput.Range("J45").Value = 50030000000
$wsInput.Range("J45").NumberFormat = "#,##0"
$wsInput.Range("J45").Font.Color = 255

$wsInput.Range("J46").Value = 590000000
$wsInput.Range("J46").Font.Color = 255

$wsInput.Range("J47").Value = 352130000
$wsInput.Range("J47").Font.Color = 255

# Restore the scroll position / active cell recorded in the sheet view.
$excel.Goto($wsInput.Range("B30"), $true)
$wsInput.Range("J50").Select()

# ---------------------------------------------------------------------
# ITR target input data
# ---------------------------------------------------------------------
$wsTarget.Activate()
$excel.Goto($wsTarget.Range("B7"), $true)
$wsTarget.Range("L14").Select()

# ---------------------------------------------------------------------
# Portfolio: RANDBETWEEN(35000,250000) sample data is volatile and is
# refreshed by the recalculation that follows this script.
# ---------------------------------------------------------------------
$wsInput.Activate()
$excel.CalculateFull()
